# "VARILLA Y ALAMBRE VISILLO" price list — monthly refresh:
#   - bump the price-list date in A1 by one month
#   - update the VARILLA CHATA price (D22)
#   - update the ALAMBRE price (D38)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 45436
$ws.Range("D22").Value = 315
$ws.Range("D38").Value = 367.127
